$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '63.551.35'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = '@'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.648.46'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = '@'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = '@'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '602.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = '@'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '146.98'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = '@'
$ws.Range("E6").Value = '  +1.71%  '
$ws.Range("E6").ClearFormats()

$ws.Range("E7").NumberFormat = '@'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E7").ClearFormats()

$ws.Range("E8").NumberFormat = '@'
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E8").ClearFormats()

$ws.Range("E9").NumberFormat = '@'
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("E9").ClearFormats()

$ws.Range("E10").NumberFormat = '@'
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("E10").ClearFormats()

$ws.Range("E11").NumberFormat = '@'
$ws.Range("E11").Value = '  +4.34%  '
$ws.Range("E11").ClearFormats()

$ws.Range("E12").NumberFormat = '@'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '27.48'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = '@'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '3.126.13'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = '@'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '63.445.70'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = '@'
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("E15").ClearFormats()

$ws.Range("E16").NumberFormat = '@'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '2.654.52'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("E17").ClearFormats()

$ws.Range("E18").NumberFormat = '@'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("E18").ClearFormats()

$ws.Range("E19").NumberFormat = '@'
$ws.Range("E19").Value = '  +4.21%  '
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '341.71'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = '@'
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E20").ClearFormats()

$ws.Range("E21").NumberFormat = '@'
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("E21").ClearFormats()

$ws.Range("E22").NumberFormat = '@'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '5.59'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = '@'
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '66.71'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = '@'
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("E24").ClearFormats()

$ws.Range("E25").NumberFormat = '@'
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '9.03'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = '@'
$ws.Range("E26").Value = '  +6.87%  '
$ws.Range("E26").ClearFormats()

$ws.Range("B27").Value = 'SuiNetwork'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '1.55'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = '@'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E27").ClearFormats()

$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '560.09'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("E28").ClearFormats()

$ws.Range("E29").NumberFormat = '@'
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("E29").ClearFormats()

$ws.Range("E30").NumberFormat = '@'
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E30").ClearFormats()

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '7.88'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = '@'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '2.02'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = '@'
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("E32").ClearFormats()

$ws.Range("E33").NumberFormat = '@'
$ws.Range("E33").Value = '  -4.05%  '
$ws.Range("E33").ClearFormats()

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.0₃0815'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = '@'
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E34").ClearFormats()

$ws.Range("E35").NumberFormat = '@'
$ws.Range("E35").Value = '  +5.19%  '
$ws.Range("E35").ClearFormats()

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '167.34'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = '@'
$ws.Range("E36").Value = '  -4.00%  '
$ws.Range("E36").ClearFormats()

$ws.Range("E37").NumberFormat = '@'
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("E37").ClearFormats()

$ws.Range("E38").NumberFormat = '@'
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.93'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = '@'
$ws.Range("E39").Value = '  +5.88%  '
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '19.10'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = '@'
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E40").ClearFormats()

$ws.Range("E41").NumberFormat = '@'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '168.41'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = '@'
$ws.Range("E42").Value = '  -1.33%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '3.76'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = '@'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '22.14'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = '@'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E44").ClearFormats()

$ws.Range("E45").NumberFormat = '@'
$ws.Range("E45").Value = '  +3.08%  '
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.630'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = '@'
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = '@'
$ws.Range("E47").Value = '  +3.45%  '
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.0960'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = '@'
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '18.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = '@'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E49").ClearFormats()

$ws.Range("E50").NumberFormat = '@'
$ws.Range("E50").Value = '  +9.45%  '
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '11.27'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = '@'
$ws.Range("E51").Value = '  -0.74%  '
$ws.Range("E51").ClearFormats()
